# Auto-generated Excel COM-interop script applying the cryptos.xlsx diff
# (Mon Jul  1 20:52:36 UTC 2024 GitHub Actions crypto price/volume update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.305.24"
$ws.Range("E2").Value = "  +2.07%  "
$ws.Range("D3").Value = "3.467.84"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.44"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.76"
$ws.Range("E6").Value = "  +2.09%  "
$ws.Range("D7").Value = "3.467.63"
$ws.Range("E7").Value = "  +1.28%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("E10").Value = "  +1.81%  "
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.406"
$ws.Range("E12").Value = "  +4.86%  "
$ws.Range("D13").Value = "4.062.59"
$ws.Range("E13").Value = "  +1.30%  "
$ws.Range("E14").Value = "  +2.76%  "
$ws.Range("E15").Value = "  +2.59%  "
$ws.Range("D16").Value = "3.463.19"
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000173"
$ws.Range("E17").Value = "  +0.98%  "
$ws.Range("D18").Value = "63.298.90"
$ws.Range("E18").Value = "  +1.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.41"
$ws.Range("E19").Value = "  +3.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.54"
$ws.Range("E20").Value = "  +3.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.35"
$ws.Range("E21").Value = "  +1.58%  "
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.567"
$ws.Range("E23").Value = "  +1.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.53"
$ws.Range("E24").Value = "  -0.50%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "3.613.66"
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000117"
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.183"
$ws.Range("E28").Value = "  -2.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.67"
$ws.Range("E29").Value = "  +1.92%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.20"
$ws.Range("E31").Value = "  +1.92%  "
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.48"
$ws.Range("E34").Value = "  -0.62%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.34"
$ws.Range("E35").Value = "  -4.35%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.34"
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.17"
$ws.Range("E37").Value = "  +2.34%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.61"
$ws.Range("E38").Value = "  +4.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "31.98"
$ws.Range("E39").Value = "  +8.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "167.37"
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("D41").Value = "3.506.36"
$ws.Range("E41").Value = "  +1.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0769"
$ws.Range("E42").Value = "  +1.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.794"
$ws.Range("E43").Value = "  +0.77%  "
$ws.Range("E44").Value = "  +4.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.41"
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("E46").Value = "  +3.31%  "
$ws.Range("E47").Value = "  -1.39%  "
$ws.Range("D48").Value = "2.591.76"
$ws.Range("E48").Value = "  +3.17%  "
$ws.Range("E49").Value = "  +9.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.84"
$ws.Range("E50").Value = "  +2.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.07"
$ws.Range("E51").Value = "  +0.09%  "
